$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''245.40'
$ws.Range("E2").Value = '''0.56%'

$ws.Range("D3").Value = '''29.85'
$ws.Range("E3").Value = '''12.87%'

$ws.Range("D4").Value = '''5.134'
$ws.Range("E4").Value = '''-0.05%'

$ws.Range("D5").Value = '''0.05711'
$ws.Range("E5").Value = '''1.88%'

$ws.Range("D6").Value = '''6.563'
$ws.Range("E6").Value = '''1.46%'

$ws.Range("D7").Value = '''0.8562'
$ws.Range("E7").Value = '''4.57%'

$ws.Range("D8").Value = '''0.8727'
$ws.Range("E8").Value = '''4.82%'

$ws.Range("D9").Value = '''0.1345'
$ws.Range("E9").Value = '''0.96%'

$ws.Range("D10").Value = '''0.06914'
$ws.Range("E10").Value = '''-0.32%'

$ws.Range("D11").Value = '''0.02892'
$ws.Range("E11").Value = '''-0.12%'

$ws.Range("D12").Value = '''0.09376'
$ws.Range("E12").Value = '''-0.03%'

$ws.Range("D13").Value = '''0.001506'
$ws.Range("E13").Value = '''-1.10%'

$ws.Range("D14").Value = '''0.04174'
$ws.Range("E14").Value = '''-9.09%'

$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.005981'
$ws.Range("E15").Value = '''-2.83%'

$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.512'
$ws.Range("E16").Value = '''-3.85%'

$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.016'
$ws.Range("E17").Value = '''-0.28%'

$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.181'
$ws.Range("E18").Value = '''-5.20%'

$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").Value = '''0.01016'
$ws.Range("E19").Value = '''1,604.43%'

$ws.Range("D21").Value = '''0.03366'
$ws.Range("E21").Value = '''8.90%'

$ws.Range("D22").Value = '''0.1304'
$ws.Range("E22").Value = '''0.36%'

$ws.Range("D23").Value = '''3.601'
$ws.Range("E23").Value = '''-3.82%'

$ws.Range("E24").Value = '''2.35%'

$ws.Range("D25").Value = '''0.001206'
$ws.Range("E25").Value = '''-1.58%'

$ws.Range("D26").Value = '''0.004480'
$ws.Range("E26").Value = '''-0.31%'

$ws.Range("D27").Value = '''0.0001176'
$ws.Range("E27").Value = '''22.47%'

$ws.Range("D28").Value = '''0.0001387'
$ws.Range("E28").Value = '''-0.90%'

$ws.Range("D40").Value = '''0.03767'
$ws.Range("E40").Value = '''3.45%'

$ws.Range("D41").Value = '''0.005799'
$ws.Range("E41").Value = '''-6.01%'

$ws.Range("D42").Value = '''0.1064'
$ws.Range("E42").Value = '''1.24%'

$ws.Range("D43").Value = '''0.002224'
$ws.Range("E43").Value = '''-7.35%'

$ws.Range("D44").Value = '''0.009549'
$ws.Range("E44").Value = '''17.69%'

$ws.Range("D45").Value = '''0.00005065'
$ws.Range("E45").Value = '''-5.38%'

$ws.Range("E46").Value = '''-0.34%'

$ws.Range("D47").Value = '''0.07973'
$ws.Range("E47").Value = '''-43.05%'

$ws.Range("E48").Value = '''11.86%'

$ws.Range("D49").Value = '''0.00002093'
$ws.Range("E49").Value = '''-0.34%'

$ws.Range("D50").Value = '''0.0001993'
$ws.Range("E50").Value = '''-0.34%'

